$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.200.93'
$ws.Range('E2').Value = '  -2.11%  '
$ws.Range('D3').Value = '3.386.47'
$ws.Range('E3').Value = '  -1.81%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'570.72"
$ws.Range('E5').Value = '  -1.55%  '
$ws.Range('D6').Value = "'141.22"
$ws.Range('E6').Value = '  -5.00%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '3.386.27'
$ws.Range('E8').Value = '  -1.85%  '
$ws.Range('D9').Value = "'0.475"
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').Value = "'7.47"
$ws.Range('E10').Value = '  -4.41%  '
$ws.Range('D11').Value = "'0.124"
$ws.Range('E11').Value = '  -0.90%  '
$ws.Range('D12').Value = "'0.394"
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('D13').Value = '3.960.00'
$ws.Range('E13').Value = '  -1.89%  '
$ws.Range('D14').Value = "'28.18"
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').Value = "'0.124"
$ws.Range('E15').Value = '  +1.05%  '
$ws.Range('D16').Value = "'0.0000171"
$ws.Range('E16').Value = '  -2.43%  '
$ws.Range('D17').Value = '3.384.16'
$ws.Range('E17').Value = '  -1.96%  '
$ws.Range('D18').Value = '60.369.37'
$ws.Range('E18').Value = '  -1.94%  '
$ws.Range('D19').Value = "'6.29"
$ws.Range('E19').Value = '  -0.87%  '
$ws.Range('D20').Value = "'14.11"
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').Value = "'9.18"
$ws.Range('E21').Value = '  -2.69%  '
$ws.Range('D22').Value = "'389.18"
$ws.Range('E22').Value = '  +0.54%  '
$ws.Range('D23').Value = "'0.562"
$ws.Range('E23').Value = '  -1.48%  '
$ws.Range('D24').Value = "'73.42"
$ws.Range('E24').Value = '  +0.91%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').Value = "'0.0000118"
$ws.Range('E26').Value = '  -4.13%  '
$ws.Range('D27').Value = '3.517.15'
$ws.Range('E27').Value = '  -2.13%  '
$ws.Range('E28').Value = '  -1.03%  '
$ws.Range('D29').Value = "'1.00"
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').Value = "'7.41"
$ws.Range('E30').Value = '  -5.12%  '
$ws.Range('D31').Value = "'8.07"
$ws.Range('E31').Value = '  -2.13%  '
$ws.Range('E32').Value = '  -1.18%  '
$ws.Range('E33').Value = '  -6.68%  '
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').Value = "'23.76"
$ws.Range('E35').Value = '  -0.92%  '
$ws.Range('D36').Value = "'6.96"
$ws.Range('E36').Value = '  -1.80%  '
$ws.Range('D37').Value = '3.411.94'
$ws.Range('E37').Value = '  -1.73%  '
$ws.Range('D38').Value = "'167.65"
$ws.Range('E38').Value = '  +0.99%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = "'4.95"
$ws.Range('E39').Value = '  -5.76%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').Value = "'1.50"
$ws.Range('E40').Value = '  -4.20%  '
$ws.Range('D41').Value = "'0.0778"
$ws.Range('E41').Value = '  -1.78%  '
$ws.Range('D42').Value = "'27.15"
$ws.Range('E42').Value = '  +4.38%  '
$ws.Range('D43').Value = "'0.784"
$ws.Range('E43').Value = '  -1.47%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').Value = "'4.47"
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('E46').Value = '  -1.93%  '
$ws.Range('D47').Value = "'41.18"
$ws.Range('E47').Value = '  -2.63%  '
$ws.Range('D48').Value = '2.529.72'
$ws.Range('E48').Value = '  -3.25%  '
$ws.Range('D49').Value = "'1.12"
$ws.Range('E49').Value = '  -3.57%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = "'23.23"
$ws.Range('E50').Value = '  +0.48%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').Value = "'6.85"
$ws.Range('E51').Value = '  -2.55%  '
